$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 5 (the "B1- TYPES OF FINANCIAL DOCUMENTS" slide): change the
#    table's style from the custom "Table_0" style to the built-in style
#    {77584228-11C4-47D9-B6B2-7107381082F8}.
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
  $shp = $slide5.Shapes.Item($i)
  if ($shp.HasTable) {
    $shp.Table.ApplyStyle("{77584228-11C4-47D9-B6B2-7107381082F8}")
  }
}

# ---------------------------------------------------------------------------
# 2. Swap the presentation's colour theme back to the stock Office palette
#    (the deck currently carries the custom "Integral / Red Violet" palette
#    on its slide master/theme; restore the default Office colours).
# ---------------------------------------------------------------------------
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $r, $g, $b) {
  $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

Set-ThemeColor $themeColors 1  0x00 0x00 0x00   # dk1
Set-ThemeColor $themeColors 2  0xFF 0xFF 0xFF   # lt1
Set-ThemeColor $themeColors 3  0x44 0x54 0x6A   # dk2
Set-ThemeColor $themeColors 4  0xE7 0xE6 0xE6   # lt2
Set-ThemeColor $themeColors 5  0x5B 0x9B 0xD5   # accent1
Set-ThemeColor $themeColors 6  0xED 0x7D 0x31   # accent2
Set-ThemeColor $themeColors 7  0xA5 0xA5 0xA5   # accent3
Set-ThemeColor $themeColors 8  0xFF 0xC0 0x00   # accent4
Set-ThemeColor $themeColors 9  0x44 0x72 0xC4   # accent5
Set-ThemeColor $themeColors 10 0x70 0xAD 0x47   # accent6
Set-ThemeColor $themeColors 11 0x05 0x63 0xC1   # hlink
Set-ThemeColor $themeColors 12 0x95 0x4F 0x72   # folHlink
